$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row labels (Holden, Rizzie Spiral inserted; Thomas Hex -> Matthies Hex rename)
$ws.Cells.Item(4, 2).Value2 = "Holden"
$ws.Cells.Item(5, 2).Value2 = "Rizzie Spiral"
$ws.Cells.Item(11, 2).Value2 = "Matthies Hex"

# New rows 30 and 31 labels + index column
$ws.Cells.Item(30, 1).Value2 = 28
$ws.Cells.Item(30, 2).Value2 = "Michael-CCHex"
$ws.Cells.Item(31, 1).Value2 = 29
$ws.Cells.Item(31, 2).Value2 = "Michael-SNHex"

# Rewrite the C:W data grid for rows 4 through 31 (values shifted after the two inserted techniques)
$ws.Cells.Item(4, 3).Value2 = 0.9938350625127605
$ws.Cells.Item(4, 4).Value2 = 1.001519876326294
$ws.Cells.Item(4, 5).Value2 = 1.004147923631342
$ws.Cells.Item(4, 6).Value2 = 0.9987533462632899
$ws.Cells.Item(4, 7).Value2 = 0.9938350625127605
$ws.Cells.Item(4, 8).Value2 = 1.001553098832244
$ws.Cells.Item(4, 9).Value2 = 0.99598782529064
$ws.Cells.Item(4, 10).Value2 = 1.004147923631342
$ws.Cells.Item(4, 11).Value2 = 1.004147923631342
$ws.Cells.Item(4, 12).Value2 = 1.002286127790005
$ws.Cells.Item(4, 13).Value2 = 0.9986740186969879
$ws.Cells.Item(4, 14).Value2 = 1.004147923631342
$ws.Cells.Item(4, 15).Value2 = 1.001519876326294
$ws.Cells.Item(4, 16).Value2 = 0.9976774694195273
$ws.Cells.Item(4, 17).Value2 = 1.000096947511641
$ws.Cells.Item(4, 18).Value2 = 0.9998342874901321
$ws.Cells.Item(4, 19).Value2 = 0.9980096525120142
$ws.Cells.Item(4, 20).Value2 = 0.9998342874901321
$ws.Cells.Item(4, 21).Value2 = 0.9995442202918461
$ws.Cells.Item(4, 22).Value2 = 1.000464960959745
$ws.Cells.Item(4, 23).Value2 = 0.9995946599179454

$ws.Cells.Item(5, 3).Value2 = 0.9835397744832932
$ws.Cells.Item(5, 4).Value2 = 1.003863953474714
$ws.Cells.Item(5, 5).Value2 = 1.011419896858945
$ws.Cells.Item(5, 6).Value2 = 0.9965472491412568
$ws.Cells.Item(5, 7).Value2 = 0.9835397744832932
$ws.Cells.Item(5, 8).Value2 = 1.004254568284014
$ws.Cells.Item(5, 9).Value2 = 0.9892539844968623
$ws.Cells.Item(5, 10).Value2 = 1.011419896858945
$ws.Cells.Item(5, 11).Value2 = 1.011419896858945
$ws.Cells.Item(5, 12).Value2 = 1.006083523288294
$ws.Cells.Item(5, 13).Value2 = 0.9964479929250627
$ws.Cells.Item(5, 14).Value2 = 1.011419896858945
$ws.Cells.Item(5, 15).Value2 = 1.003863953474714
$ws.Cells.Item(5, 16).Value2 = 0.9937018639790038
$ws.Cells.Item(5, 17).Value2 = 1.000155973199889
$ws.Cells.Item(5, 18).Value2 = 0.9996078749389842
$ws.Cells.Item(5, 19).Value2 = 0.9946172402943567
$ws.Cells.Item(5, 20).Value2 = 0.9996078749389842
$ws.Cells.Item(5, 21).Value2 = 0.9988179044355039
$ws.Cells.Item(5, 22).Value2 = 1.001338302920192
$ws.Cells.Item(5, 23).Value2 = 0.9989263678690554

$ws.Cells.Item(6, 3).Value2 = 0.9925287322123724
$ws.Cells.Item(6, 4).Value2 = 1.001944550056633
$ws.Cells.Item(6, 5).Value2 = 1.004844426022796
$ws.Cells.Item(6, 6).Value2 = 0.998554859342634
$ws.Cells.Item(6, 7).Value2 = 0.9925287322123724
$ws.Cells.Item(6, 8).Value2 = 1.001825183956162
$ws.Cells.Item(6, 9).Value2 = 0.9951554259965935
$ws.Cells.Item(6, 10).Value2 = 1.004844426022796
$ws.Cells.Item(6, 11).Value2 = 1.004844426022796
$ws.Cells.Item(6, 12).Value2 = 1.002781319427469
$ws.Cells.Item(6, 13).Value2 = 0.9983992164839069
$ws.Cells.Item(6, 14).Value2 = 1.004844426022796
$ws.Cells.Item(6, 15).Value2 = 1.001944550056633
$ws.Cells.Item(6, 16).Value2 = 0.9972366411345025
$ws.Cells.Item(6, 17).Value2 = 1.00017188327027
$ws.Cells.Item(6, 18).Value2 = 0.9997725694306002
$ws.Cells.Item(6, 19).Value2 = 0.9976241662509707
$ws.Cells.Item(6, 20).Value2 = 0.9997725694306002
$ws.Cells.Item(6, 21).Value2 = 0.9994292311939269
$ws.Cells.Item(6, 22).Value2 = 1.000512270159701
$ws.Cells.Item(6, 23).Value2 = 0.9995042141873207

$ws.Cells.Item(7, 3).Value2 = 0.9918932464985607
$ws.Cells.Item(7, 4).Value2 = 1.001995315367437
$ws.Cells.Item(7, 5).Value2 = 1.005460269351582
$ws.Cells.Item(7, 6).Value2 = 0.9983585744092227
$ws.Cells.Item(7, 7).Value2 = 0.9918932464985607
$ws.Cells.Item(7, 8).Value2 = 1.002044117074929
$ws.Cells.Item(7, 9).Value2 = 0.9947235115850159
$ws.Cells.Item(7, 10).Value2 = 1.005460269351582
$ws.Cells.Item(7, 11).Value2 = 1.005460269351582
$ws.Cells.Item(7, 12).Value2 = 1.003005860453888
$ws.Cells.Item(7, 13).Value2 = 0.9982561669668547
$ws.Cells.Item(7, 14).Value2 = 1.005460269351582
$ws.Cells.Item(7, 15).Value2 = 1.001995315367437
$ws.Cells.Item(7, 16).Value2 = 0.9969442809329987
$ws.Cells.Item(7, 17).Value2 = 1.000125741167146
$ws.Cells.Item(7, 18).Value2 = 0.9997829437391932
$ws.Cells.Item(7, 19).Value2 = 0.9973815762776174
$ws.Cells.Item(7, 20).Value2 = 0.9997829437391932
$ws.Cells.Item(7, 21).Value2 = 0.9994012495461085
$ws.Cells.Item(7, 22).Value2 = 1.000613053507203
$ws.Cells.Item(7, 23).Value2 = 0.9994671327134362

$ws.Cells.Item(8, 3).Value2 = 0.9734344280050384
$ws.Cells.Item(8, 4).Value2 = 1.006502701740315
$ws.Cells.Item(8, 5).Value2 = 1.017956891462178
$ws.Cells.Item(8, 6).Value2 = 0.9945981417537488
$ws.Cells.Item(8, 7).Value2 = 0.9734344280050384
$ws.Cells.Item(8, 8).Value2 = 1.006718450423402
$ws.Cells.Item(8, 9).Value2 = 0.9827029036667746
$ws.Cells.Item(8, 10).Value2 = 1.017956891462178
$ws.Cells.Item(8, 11).Value2 = 1.017956891462178
$ws.Cells.Item(8, 12).Value2 = 1.009846329914978
$ws.Cells.Item(8, 13).Value2 = 0.994283358648875
$ws.Cells.Item(8, 14).Value2 = 1.017956891462178
$ws.Cells.Item(8, 15).Value2 = 1.006502701740315
$ws.Cells.Item(8, 16).Value2 = 0.9899685648726769
$ws.Cells.Item(8, 17).Value2 = 1.000393030194595
$ws.Cells.Item(8, 18).Value2 = 0.9992980070691772
$ws.Cells.Item(8, 19).Value2 = 0.991406829464743
$ws.Cells.Item(8, 20).Value2 = 0.9992980070691772
$ws.Cells.Item(8, 21).Value2 = 0.9980443449641015
$ws.Cells.Item(8, 22).Value2 = 1.002026854263717
$ws.Cells.Item(8, 23).Value2 = 0.9982554007019138

$ws.Cells.Item(9, 3).Value2 = 0.9991283159091601
$ws.Cells.Item(9, 4).Value2 = 1.000218842449253
$ws.Cells.Item(9, 5).Value2 = 1.00057948312
$ws.Cells.Item(9, 6).Value2 = 0.999826255486405
$ws.Cells.Item(9, 7).Value2 = 0.9991283159091601
$ws.Cells.Item(9, 8).Value2 = 1.000217408302308
$ws.Cells.Item(9, 9).Value2 = 0.9994333865012439
$ws.Cells.Item(9, 10).Value2 = 1.00057948312
$ws.Cells.Item(9, 11).Value2 = 1.00057948312
$ws.Cells.Item(9, 12).Value2 = 1.000323657782904
$ws.Cells.Item(9, 13).Value2 = 0.9998127518177365
$ws.Cells.Item(9, 14).Value2 = 1.00057948312
$ws.Cells.Item(9, 15).Value2 = 1.000218842449253
$ws.Cells.Item(9, 16).Value2 = 0.9996735791792064
$ws.Cells.Item(9, 17).Value2 = 1.000015797133495
$ws.Cells.Item(9, 18).Value2 = 0.999975547159471
$ws.Cells.Item(9, 19).Value2 = 0.9997199700587164
$ws.Cells.Item(9, 20).Value2 = 0.999975547159471
$ws.Cells.Item(9, 21).Value2 = 0.9999348483240373
$ws.Cells.Item(9, 22).Value2 = 1.00006377528323
$ws.Cells.Item(9, 23).Value2 = 0.9999425126711263

$ws.Cells.Item(10, 3).Value2 = 0.9999353741012039
$ws.Cells.Item(10, 4).Value2 = 1.000013664376807
$ws.Cells.Item(10, 5).Value2 = 1.000047512949314
$ws.Cells.Item(10, 6).Value2 = 0.9999854820288959
$ws.Cells.Item(10, 7).Value2 = 0.9999353741012039
$ws.Cells.Item(10, 8).Value2 = 1.000017540400062
$ws.Cells.Item(10, 9).Value2 = 0.9999575497064418
$ws.Cells.Item(10, 10).Value2 = 1.000047512949314
$ws.Cells.Item(10, 11).Value2 = 1.000047512949314
$ws.Cells.Item(10, 12).Value2 = 1.000023727047024
$ws.Cells.Item(10, 13).Value2 = 0.999985963712824
$ws.Cells.Item(10, 14).Value2 = 1.000047512949314
$ws.Cells.Item(10, 15).Value2 = 1.000013664376807
$ws.Cells.Item(10, 16).Value2 = 0.9999745192390053
$ws.Cells.Item(10, 17).Value2 = 0.9999998140448154
$ws.Cells.Item(10, 18).Value2 = 0.9999988504757749
$ws.Cells.Item(10, 19).Value2 = 0.9999783340636115
$ws.Cells.Item(10, 20).Value2 = 0.9999988504757749
$ws.Cells.Item(10, 21).Value2 = 0.9999956287850371
$ws.Cells.Item(10, 22).Value2 = 1.000006005617892
$ws.Cells.Item(10, 23).Value2 = 0.9999958517903216

$ws.Cells.Item(11, 3).Value2 = 0.9985294072045284
$ws.Cells.Item(11, 4).Value2 = 1.000373151298786
$ws.Cells.Item(11, 5).Value2 = 1.000970607713255
$ws.Cells.Item(11, 6).Value2 = 0.999709407598304
$ws.Cells.Item(11, 7).Value2 = 0.9985294072045284
$ws.Cells.Item(11, 8).Value2 = 1.00036458846483
$ws.Cells.Item(11, 9).Value2 = 0.9990447674872985
$ws.Cells.Item(11, 10).Value2 = 1.000970607713255
$ws.Cells.Item(11, 11).Value2 = 1.000970607713255
$ws.Cells.Item(11, 12).Value2 = 1.0005464477392
$ws.Cells.Item(11, 13).Value2 = 0.9996843367449811
$ws.Cells.Item(11, 14).Value2 = 1.000970607713255
$ws.Cells.Item(11, 15).Value2 = 1.000373151298786
$ws.Cells.Item(11, 16).Value2 = 0.9994512792516574
$ws.Cells.Item(11, 17).Value2 = 1.000028744021884
$ws.Cells.Item(11, 18).Value2 = 0.9999577220721901
$ws.Cells.Item(11, 19).Value2 = 0.9995289650827653
$ws.Cells.Item(11, 20).Value2 = 0.9999577220721901
$ws.Cells.Item(11, 21).Value2 = 0.9998893757403878
$ws.Cells.Item(11, 22).Value2 = 1.000105622134961
$ws.Cells.Item(11, 23).Value2 = 0.999902839281398

$ws.Cells.Item(12, 3).Value2 = 0.9730416487440269
$ws.Cells.Item(12, 4).Value2 = 1.00659657855665
$ws.Cells.Item(12, 5).Value2 = 1.018226420655612
$ws.Cells.Item(12, 6).Value2 = 0.9945168220059994
$ws.Cells.Item(12, 7).Value2 = 0.9730416487440269
$ws.Cells.Item(12, 8).Value2 = 1.006819044154509
$ws.Cells.Item(12, 9).Value2 = 0.9824467685433133
$ws.Cells.Item(12, 10).Value2 = 1.018226420655612
$ws.Cells.Item(12, 11).Value2 = 1.018226420655612
$ws.Cells.Item(12, 12).Value2 = 1.009991673261425
$ws.Cells.Item(12, 13).Value2 = 0.9941986996291773
$ws.Cells.Item(12, 14).Value2 = 1.018226420655612
$ws.Cells.Item(12, 15).Value2 = 1.00659657855665
$ws.Cells.Item(12, 16).Value2 = 0.9898191136503387
$ws.Cells.Item(12, 17).Value2 = 1.000397639092914
$ws.Cells.Item(12, 18).Value2 = 0.9992882159854298
$ws.Cells.Item(12, 19).Value2 = 0.9912789756432848
$ws.Cells.Item(12, 20).Value2 = 0.9992882159854298
$ws.Cells.Item(12, 21).Value2 = 0.9980158368963666
$ws.Cells.Item(12, 22).Value2 = 1.002057953648216
$ws.Cells.Item(12, 23).Value2 = 0.9982297069438393

$ws.Cells.Item(13, 3).Value2 = 1.004255291846079
$ws.Cells.Item(13, 4).Value2 = 0.9989910918505308
$ws.Cells.Item(13, 5).Value2 = 0.9970655070986497
$ws.Cells.Item(13, 6).Value2 = 1.000886207663159
$ws.Cells.Item(13, 7).Value2 = 1.004255291846079
$ws.Cells.Item(13, 8).Value2 = 0.9989056650958622
$ws.Cells.Item(13, 9).Value2 = 1.002776325247888
$ws.Cells.Item(13, 10).Value2 = 0.9970655070986497
$ws.Cells.Item(13, 11).Value2 = 0.9970655070986497
$ws.Cells.Item(13, 12).Value2 = 0.9984262332657451
$ws.Cells.Item(13, 13).Value2 = 1.000917662219304
$ws.Cells.Item(13, 14).Value2 = 0.9970655070986497
$ws.Cells.Item(13, 15).Value2 = 0.9989910918505308
$ws.Cells.Item(13, 16).Value2 = 1.001623191848305
$ws.Cells.Item(13, 17).Value2 = 0.9999543770349175
$ws.Cells.Item(13, 18).Value2 = 1.00010396359842
$ws.Cells.Item(13, 19).Value2 = 1.001388015305305
$ws.Cells.Item(13, 20).Value2 = 1.00010396359842
$ws.Cells.Item(13, 21).Value2 = 1.000307388253641
$ws.Cells.Item(13, 22).Value2 = 0.9996590120226425
$ws.Cells.Item(13, 23).Value2 = 1.000277998035902

$ws.Cells.Item(14, 3).Value2 = 0.9921544834210541
$ws.Cells.Item(14, 4).Value2 = 1.001915478452633
$ws.Cells.Item(14, 5).Value2 = 1.005311937852631
$ws.Cells.Item(14, 6).Value2 = 0.9984015232842081
$ws.Cells.Item(14, 7).Value2 = 0.9921544834210541
$ws.Cells.Item(14, 8).Value2 = 1.001986879284209
$ws.Cells.Item(14, 9).Value2 = 0.994890853452632
$ws.Cells.Item(14, 10).Value2 = 1.005311937852631
$ws.Cells.Item(14, 11).Value2 = 1.005311937852631
$ws.Cells.Item(14, 12).Value2 = 1.00290736758947
$ws.Cells.Item(14, 13).Value2 = 0.998311426305264
$ws.Cells.Item(14, 14).Value2 = 1.005311937852631
$ws.Cells.Item(14, 15).Value2 = 1.001915478452633
$ws.Cells.Item(14, 16).Value2 = 0.9970349809368437
$ws.Cells.Item(14, 17).Value2 = 1.000113452378949
$ws.Cells.Item(14, 18).Value2 = 0.9997939665754393
$ws.Cells.Item(14, 19).Value2 = 0.9974604627263172
$ws.Cells.Item(14, 20).Value2 = 0.9997939665754393
$ws.Cells.Item(14, 21).Value2 = 0.9994233315078955
$ws.Cells.Item(14, 22).Value2 = 1.000601052776843
$ws.Cells.Item(14, 23).Value2 = 0.9994849937052626

$ws.Cells.Item(15, 3).Value2 = 1.00531680606041
$ws.Cells.Item(15, 4).Value2 = 0.9986786091333335
$ws.Cells.Item(15, 5).Value2 = 0.9964415780368485
$ws.Cells.Item(15, 6).Value2 = 1.001068359872898
$ws.Cells.Item(15, 7).Value2 = 1.00531680606041
$ws.Cells.Item(15, 8).Value2 = 0.9986664585994167
$ws.Cells.Item(15, 9).Value2 = 1.003458370875478
$ws.Cells.Item(15, 10).Value2 = 0.9964415780368485
$ws.Cells.Item(15, 11).Value2 = 0.9964415780368485
$ws.Cells.Item(15, 12).Value2 = 0.9980272709185766
$ws.Cells.Item(15, 13).Value2 = 1.00114292217613
$ws.Cells.Item(15, 14).Value2 = 0.9964415780368485
$ws.Cells.Item(15, 15).Value2 = 0.9986786091333335
$ws.Cells.Item(15, 16).Value2 = 1.001997707596872
$ws.Cells.Item(15, 17).Value2 = 0.9999107656547319
$ws.Cells.Item(15, 18).Value2 = 1.000145664410197
$ws.Cells.Item(15, 19).Value2 = 1.001712779123291
$ws.Cells.Item(15, 20).Value2 = 1.000145664410197
$ws.Cells.Item(15, 21).Value2 = 1.000394978851681
$ws.Cells.Item(15, 22).Value2 = 0.9996042986887141
$ws.Cells.Item(15, 23).Value2 = 1.000350046959137

$ws.Cells.Item(16, 3).Value2 = 0.9534276099999992
$ws.Cells.Item(16, 4).Value2 = 1.011390400000001
$ws.Cells.Item(16, 5).Value2 = 1.0314974
$ws.Cells.Item(16, 6).Value2 = 0.9905238300000004
$ws.Cells.Item(16, 7).Value2 = 0.9534276099999992
$ws.Cells.Item(16, 8).Value2 = 1.011783500000001
$ws.Cells.Item(16, 9).Value2 = 0.9696746199999998
$ws.Cells.Item(16, 10).Value2 = 1.0314974
$ws.Cells.Item(16, 11).Value2 = 1.0314974
$ws.Cells.Item(16, 12).Value2 = 1.0172607
$ws.Cells.Item(16, 13).Value2 = 0.9899775199999981
$ws.Cells.Item(16, 14).Value2 = 1.0314974
$ws.Cells.Item(16, 15).Value2 = 1.011390400000001
$ws.Cells.Item(16, 16).Value2 = 0.9824090050000001
$ws.Cells.Item(16, 17).Value2 = 1.000683959999999
$ws.Cells.Item(16, 18).Value2 = 0.9987718033333334
$ws.Cells.Item(16, 19).Value2 = 0.9849318433333328
$ws.Cells.Item(16, 20).Value2 = 0.9987718033333334
$ws.Cells.Item(16, 21).Value2 = 0.9965732324999996
$ws.Cells.Item(16, 22).Value2 = 1.003558066
$ws.Cells.Item(16, 23).Value2 = 0.9969419474999999

$ws.Cells.Item(17, 3).Value2 = 0.8759446600000002
$ws.Cells.Item(17, 4).Value2 = 1.0301398
$ws.Cells.Item(17, 5).Value2 = 1.0842575
$ws.Cells.Item(17, 6).Value2 = 0.97462952
$ws.Cells.Item(17, 7).Value2 = 0.8759446600000002
$ws.Cells.Item(17, 8).Value2 = 1.0314994
$ws.Cells.Item(17, 9).Value2 = 0.91918716
$ws.Cells.Item(17, 10).Value2 = 1.0842575
$ws.Cells.Item(17, 11).Value2 = 1.0842575
$ws.Cells.Item(17, 12).Value2 = 1.0459565
$ws.Cells.Item(17, 13).Value2 = 0.9732909399999999
$ws.Cells.Item(17, 14).Value2 = 1.0842575
$ws.Cells.Item(17, 15).Value2 = 1.0301398
$ws.Cells.Item(17, 16).Value2 = 0.9530422300000001
$ws.Cells.Item(17, 17).Value2 = 1.00171537
$ws.Cells.Item(17, 18).Value2 = 0.9967806533333334
$ws.Cells.Item(17, 19).Value2 = 0.9597918000000001
$ws.Cells.Item(17, 20).Value2 = 0.9967806533333334
$ws.Cells.Item(17, 21).Value2 = 0.9909082250000001
$ws.Cells.Item(17, 22).Value2 = 1.00957808
$ws.Cells.Item(17, 23).Value2 = 0.9918631849999999

$ws.Cells.Item(18, 3).Value2 = 0.9526597999999999
$ws.Cells.Item(18, 4).Value2 = 1.012634
$ws.Cells.Item(18, 5).Value2 = 1.0301398
$ws.Cells.Item(18, 6).Value2 = 0.9910432999999998
$ws.Cells.Item(18, 7).Value2 = 0.9526597999999999
$ws.Cells.Item(18, 8).Value2 = 1.0113912
$ws.Cells.Item(18, 9).Value2 = 0.9693574700000001
$ws.Cells.Item(18, 10).Value2 = 1.0301398
$ws.Cells.Item(18, 11).Value2 = 1.0301398
$ws.Cells.Item(18, 12).Value2 = 1.0176561
$ws.Cells.Item(18, 13).Value2 = 0.9898757500000001
$ws.Cells.Item(18, 14).Value2 = 1.0301398
$ws.Cells.Item(18, 15).Value2 = 1.012634
$ws.Cells.Item(18, 16).Value2 = 0.9826469
$ws.Cells.Item(18, 17).Value2 = 1.001254875
$ws.Cells.Item(18, 18).Value2 = 0.9984778666666667
$ws.Cells.Item(18, 19).Value2 = 0.9850565166666666
$ws.Cells.Item(18, 20).Value2 = 0.9984778666666667
$ws.Cells.Item(18, 21).Value2 = 0.9963273375
$ws.Cells.Item(18, 22).Value2 = 1.00308983
$ws.Cells.Item(18, 23).Value2 = 0.9968446774999999

$ws.Cells.Item(19, 3).Value2 = 0.95342312
$ws.Cells.Item(19, 4).Value2 = 1.0113915
$ws.Cells.Item(19, 5).Value2 = 1.0315004
$ws.Cells.Item(19, 6).Value2 = 0.99052292
$ws.Cells.Item(19, 7).Value2 = 0.95342312
$ws.Cells.Item(19, 8).Value2 = 1.0117846
$ws.Cells.Item(19, 9).Value2 = 0.9696716999999999
$ws.Cells.Item(19, 10).Value2 = 1.0315004
$ws.Cells.Item(19, 11).Value2 = 1.0315004
$ws.Cells.Item(19, 12).Value2 = 1.0172624
$ws.Cells.Item(19, 13).Value2 = 0.9899765600000001
$ws.Cells.Item(19, 14).Value2 = 1.0315004
$ws.Cells.Item(19, 15).Value2 = 1.0113915
$ws.Cells.Item(19, 16).Value2 = 0.9824073099999999
$ws.Cells.Item(19, 17).Value2 = 1.00068403
$ws.Cells.Item(19, 18).Value2 = 0.9987716733333333
$ws.Cells.Item(19, 19).Value2 = 0.9849303933333333
$ws.Cells.Item(19, 20).Value2 = 0.9987716733333333
$ws.Cells.Item(19, 21).Value2 = 0.996572895
$ws.Cells.Item(19, 22).Value2 = 1.003558396
$ws.Cells.Item(19, 23).Value2 = 0.9969416499999999

$ws.Cells.Item(20, 3).Value2 = 0.9809506079452054
$ws.Cells.Item(20, 4).Value2 = 1.004828233972603
$ws.Cells.Item(20, 5).Value2 = 1.012582444109589
$ws.Cells.Item(20, 6).Value2 = 0.996232296986301
$ws.Cells.Item(20, 7).Value2 = 0.9809506079452054
$ws.Cells.Item(20, 8).Value2 = 1.00472573369863
$ws.Cells.Item(20, 9).Value2 = 0.9876253915068492
$ws.Cells.Item(20, 10).Value2 = 1.012582444109589
$ws.Cells.Item(20, 11).Value2 = 1.012582444109589
$ws.Cells.Item(20, 12).Value2 = 1.007077878630137
$ws.Cells.Item(20, 13).Value2 = 0.9959107046575345
$ws.Cells.Item(20, 14).Value2 = 1.012582444109589
$ws.Cells.Item(20, 15).Value2 = 1.004828233972603
$ws.Cells.Item(20, 16).Value2 = 0.9928894209589041
$ws.Cells.Item(20, 17).Value2 = 1.000369469315069
$ws.Cells.Item(20, 18).Value2 = 0.9994537620091323
$ws.Cells.Item(20, 19).Value2 = 0.9938965155251142
$ws.Cells.Item(20, 20).Value2 = 0.9994537620091323
$ws.Cells.Item(20, 21).Value2 = 0.9985679976712328
$ws.Cells.Item(20, 22).Value2 = 1.001370886958904
$ws.Cells.Item(20, 23).Value2 = 0.9987416614383561

$ws.Cells.Item(21, 3).Value2 = 0.9530422889473683
$ws.Cells.Item(21, 4).Value2 = 1.01201257368421
$ws.Cells.Item(21, 5).Value2 = 1.030819547368421
$ws.Cells.Item(21, 6).Value2 = 0.9907832736842105
$ws.Cells.Item(21, 7).Value2 = 0.9530422889473683
$ws.Cells.Item(21, 8).Value2 = 1.011587684210526
$ws.Cells.Item(21, 9).Value2 = 0.9695151210526316
$ws.Cells.Item(21, 10).Value2 = 1.030819547368421
$ws.Cells.Item(21, 11).Value2 = 1.030819547368421
$ws.Cells.Item(21, 12).Value2 = 1.017458952631579
$ws.Cells.Item(21, 13).Value2 = 0.989926330526316
$ws.Cells.Item(21, 14).Value2 = 1.030819547368421
$ws.Cells.Item(21, 15).Value2 = 1.01201257368421
$ws.Cells.Item(21, 16).Value2 = 0.9825274313157895
$ws.Cells.Item(21, 17).Value2 = 1.000969452105263
$ws.Cells.Item(21, 18).Value2 = 0.9986248033333333
$ws.Cells.Item(21, 19).Value2 = 0.9849937310526317
$ws.Cells.Item(21, 20).Value2 = 0.9986248033333333
$ws.Cells.Item(21, 21).Value2 = 0.996450185131579
$ws.Cells.Item(21, 22).Value2 = 1.003324057578947
$ws.Cells.Item(21, 23).Value2 = 0.9968932215131578

$ws.Cells.Item(22, 3).Value2 = 0.9790700236842105
$ws.Cells.Item(22, 4).Value2 = 1.004966313157895
$ws.Cells.Item(22, 5).Value2 = 1.014426531578947
$ws.Cells.Item(22, 6).Value2 = 0.99564365
$ws.Cells.Item(22, 7).Value2 = 0.9790700236842105
$ws.Cells.Item(22, 8).Value2 = 1.005380383157895
$ws.Cells.Item(22, 9).Value2 = 0.9863451168421055
$ws.Cells.Item(22, 10).Value2 = 1.014426531578947
$ws.Cells.Item(22, 11).Value2 = 1.014426531578947
$ws.Cells.Item(22, 12).Value2 = 1.007741071052632
$ws.Cells.Item(22, 13).Value2 = 0.9954866447368419
$ws.Cells.Item(22, 14).Value2 = 1.014426531578947
$ws.Cells.Item(22, 15).Value2 = 1.004966313157895
$ws.Cells.Item(22, 16).Value2 = 0.9920181684210526
$ws.Cells.Item(22, 17).Value2 = 1.000226478947368
$ws.Cells.Item(22, 18).Value2 = 0.9994876228070174
$ws.Cells.Item(22, 19).Value2 = 0.9931743271929824
$ws.Cells.Item(22, 20).Value2 = 0.9994876228070174
$ws.Cells.Item(22, 21).Value2 = 0.9984873782894735
$ws.Cells.Item(22, 22).Value2 = 1.001675208947368
$ws.Cells.Item(22, 23).Value2 = 0.9986324667763158

$ws.Cells.Item(23, 3).Value2 = 0.9631303094325089
$ws.Cells.Item(23, 4).Value2 = 1.009630860170504
$ws.Cells.Item(23, 5).Value2 = 1.02384479496947
$ws.Cells.Item(23, 6).Value2 = 0.9928906591822821
$ws.Cells.Item(23, 7).Value2 = 0.9631303094325089
$ws.Cells.Item(23, 8).Value2 = 1.008987733906296
$ws.Cells.Item(23, 9).Value2 = 0.9760986970092749
$ws.Cells.Item(23, 10).Value2 = 1.02384479496947
$ws.Cells.Item(23, 11).Value2 = 1.02384479496947
$ws.Cells.Item(23, 12).Value2 = 1.013729099860233
$ws.Cells.Item(23, 13).Value2 = 0.9921024356548199
$ws.Cells.Item(23, 14).Value2 = 1.02384479496947
$ws.Cells.Item(23, 15).Value2 = 1.009630860170504
$ws.Cells.Item(23, 16).Value2 = 0.9863805848015066
$ws.Cells.Item(23, 17).Value2 = 1.000866647912662
$ws.Cells.Item(23, 18).Value2 = 0.9988686548574947
$ws.Cells.Item(23, 19).Value2 = 0.9882878684192776
$ws.Cells.Item(23, 20).Value2 = 0.9988686548574944
$ws.Cells.Item(23, 21).Value2 = 0.9971771000568258
$ws.Cells.Item(23, 22).Value2 = 1.002510639039355
$ws.Cells.Item(23, 23).Value2 = 0.9975518237731735

$ws.Cells.Item(24, 3).Value2 = 0.9903291887809755
$ws.Cells.Item(24, 4).Value2 = 1.002246693879998
$ws.Cells.Item(24, 5).Value2 = 1.006751230633846
$ws.Cells.Item(24, 6).Value2 = 0.9979563900521016
$ws.Cells.Item(24, 7).Value2 = 0.9903291887809755
$ws.Cells.Item(24, 8).Value2 = 1.002512707982678
$ws.Cells.Item(24, 9).Value2 = 0.9936823756656112
$ws.Cells.Item(24, 10).Value2 = 1.006751230633846
$ws.Cells.Item(24, 11).Value2 = 1.006751230633846
$ws.Cells.Item(24, 12).Value2 = 1.003571771507939
$ws.Cells.Item(24, 13).Value2 = 0.9979116945411328
$ws.Cells.Item(24, 14).Value2 = 1.006751230633846
$ws.Cells.Item(24, 15).Value2 = 1.002246693879998
$ws.Cells.Item(24, 16).Value2 = 0.9962879413304866
$ws.Cells.Item(24, 17).Value2 = 1.000079194210565
$ws.Cells.Item(24, 18).Value2 = 0.9997757044316063
$ws.Cells.Item(24, 19).Value2 = 0.996829192400702
$ws.Cells.Item(24, 20).Value2 = 0.9997757044316063
$ws.Cells.Item(24, 21).Value2 = 0.9993097019589879
$ws.Cells.Item(24, 22).Value2 = 1.00079800769396
$ws.Cells.Item(24, 23).Value2 = 0.9993702566305352

$ws.Cells.Item(25, 3).Value2 = 1.011900432336229
$ws.Cells.Item(25, 4).Value2 = 0.9970374114771908
$ws.Cells.Item(25, 5).Value2 = 0.9920441157028639
$ws.Cells.Item(25, 6).Value2 = 1.002388097741924
$ws.Cells.Item(25, 7).Value2 = 1.011900432336229
$ws.Cells.Item(25, 8).Value2 = 0.9970179290362796
$ws.Cells.Item(25, 9).Value2 = 1.007739902697986
$ws.Cells.Item(25, 10).Value2 = 0.9920441157028639
$ws.Cells.Item(25, 11).Value2 = 0.9920441157028639
$ws.Cells.Item(25, 12).Value2 = 0.9955839789808788
$ws.Cells.Item(25, 13).Value2 = 1.002557863807517
$ws.Cells.Item(25, 14).Value2 = 0.9920441157028639
$ws.Cells.Item(25, 15).Value2 = 0.9970374114771908
$ws.Cells.Item(25, 16).Value2 = 1.00446892190671
$ws.Cells.Item(25, 17).Value2 = 0.9997976376423541
$ws.Cells.Item(25, 18).Value2 = 1.000327319838761
$ws.Cells.Item(25, 19).Value2 = 1.003831902540312
$ws.Cells.Item(25, 20).Value2 = 1.000327319838761
$ws.Cells.Item(25, 21).Value2 = 1.00088495583095
$ws.Cells.Item(25, 22).Value2 = 0.9991167878053329
$ws.Cells.Item(25, 23).Value2 = 1.000783716472609

$ws.Cells.Item(26, 3).Value2 = 1.003092392309799
$ws.Cells.Item(26, 4).Value2 = 0.9992211398179152
$ws.Cells.Item(26, 5).Value2 = 0.9979486510136296
$ws.Cells.Item(26, 6).Value2 = 1.00061479192147
$ws.Cells.Item(26, 7).Value2 = 1.003092392309799
$ws.Cells.Item(26, 8).Value2 = 0.9992301057684592
$ws.Cells.Item(26, 9).Value2 = 1.002009691371328
$ws.Cells.Item(26, 10).Value2 = 0.9979486510136296
$ws.Cells.Item(26, 11).Value2 = 0.9979486510136296
$ws.Cells.Item(26, 12).Value2 = 0.9988515272108537
$ws.Cells.Item(26, 13).Value2 = 1.000664135200781
$ws.Cells.Item(26, 14).Value2 = 0.9979486510136296
$ws.Cells.Item(26, 15).Value2 = 0.9992211398179152
$ws.Cells.Item(26, 16).Value2 = 1.001156766063857
$ws.Cells.Item(26, 17).Value2 = 0.9999426375093483
$ws.Cells.Item(26, 18).Value2 = 1.000087394380448
$ws.Cells.Item(26, 19).Value2 = 1.000992555776165
$ws.Cells.Item(26, 20).Value2 = 1.000087394380448
$ws.Cells.Item(26, 21).Value2 = 1.000231579585531
$ws.Cells.Item(26, 22).Value2 = 0.9997749938711511
$ws.Cells.Item(26, 23).Value2 = 1.00020405432678

$ws.Cells.Item(27, 3).Value2 = 1.037408405204385
$ws.Cells.Item(27, 4).Value2 = 0.9904425353262464
$ws.Cells.Item(27, 5).Value2 = 0.9754261572522006
$ws.Cells.Item(27, 6).Value2 = 1.007350248002653
$ws.Cells.Item(27, 7).Value2 = 1.037408405204385
$ws.Cells.Item(27, 8).Value2 = 0.9907619887829876
$ws.Cells.Item(27, 9).Value2 = 1.024287595839436
$ws.Cells.Item(27, 10).Value2 = 0.9754261572522006
$ws.Cells.Item(27, 11).Value2 = 0.9754261572522006
$ws.Cells.Item(27, 12).Value2 = 0.9860927760842482
$ws.Cells.Item(27, 13).Value2 = 1.008025819181727
$ws.Cells.Item(27, 14).Value2 = 0.9754261572522006
$ws.Cells.Item(27, 15).Value2 = 0.9904425353262464
$ws.Cells.Item(27, 16).Value2 = 1.013925470265316
$ws.Cells.Item(27, 17).Value2 = 0.9992341772539866
$ws.Cells.Item(27, 18).Value2 = 1.001092365927611
$ws.Cells.Item(27, 19).Value2 = 1.01195891990412
$ws.Cells.Item(27, 20).Value2 = 1.001092365927611
$ws.Cells.Item(27, 21).Value2 = 1.00282572924114
$ws.Cells.Item(27, 22).Value2 = 0.9973458148433518
$ws.Cells.Item(27, 23).Value2 = 1.002474440709236

$ws.Cells.Item(28, 3).Value2 = 1.009886374988196
$ws.Cells.Item(28, 4).Value2 = 0.9976114016248713
$ws.Cells.Item(28, 5).Value2 = 0.9932615448418654
$ws.Cells.Item(28, 6).Value2 = 1.002030387587287
$ws.Cells.Item(28, 7).Value2 = 1.009886374988196
$ws.Cells.Item(28, 8).Value2 = 0.9974822932536055
$ws.Cells.Item(28, 9).Value2 = 1.006442551086529
$ws.Cells.Item(28, 10).Value2 = 0.9932615448418654
$ws.Cells.Item(28, 11).Value2 = 0.9932615448418654
$ws.Cells.Item(28, 12).Value2 = 0.9963389829089473
$ws.Cells.Item(28, 13).Value2 = 1.002129334132218
$ws.Cells.Item(28, 14).Value2 = 0.9932615448418654
$ws.Cells.Item(28, 15).Value2 = 0.9976114016248713
$ws.Cells.Item(28, 16).Value2 = 1.003748888306534
$ws.Cells.Item(28, 17).Value2 = 0.9998703678785448
$ws.Cells.Item(28, 18).Value2 = 1.000253107151644
$ws.Cells.Item(28, 19).Value2 = 1.003209036915095
$ws.Cells.Item(28, 20).Value2 = 1.000253107151644
$ws.Cells.Item(28, 21).Value2 = 1.000722163896788
$ws.Cells.Item(28, 22).Value2 = 0.9992300400858033
$ws.Cells.Item(28, 23).Value2 = 1.00064785880294

$ws.Cells.Item(29, 3).Value2 = 0.9972091470990151
$ws.Cells.Item(29, 4).Value2 = 1.000703463735362
$ws.Cells.Item(29, 5).Value2 = 1.00185034238029
$ws.Cells.Item(29, 6).Value2 = 0.9994455036034408
$ws.Cells.Item(29, 7).Value2 = 0.9972091470990151
$ws.Cells.Item(29, 8).Value2 = 1.000694522585504
$ws.Cells.Item(29, 9).Value2 = 0.9981863705629614
$ws.Cells.Item(29, 10).Value2 = 1.00185034238029
$ws.Cells.Item(29, 11).Value2 = 1.00185034238029
$ws.Cells.Item(29, 12).Value2 = 1.001036534569038
$ws.Cells.Item(29, 13).Value2 = 0.9994006546502689
$ws.Cells.Item(29, 14).Value2 = 1.00185034238029
$ws.Cells.Item(29, 15).Value2 = 1.000703463735362
$ws.Cells.Item(29, 16).Value2 = 0.9989563054171888
$ws.Cells.Item(29, 17).Value2 = 1.000052059192816
$ws.Cells.Item(29, 18).Value2 = 0.9999209844048892
$ws.Cells.Item(29, 19).Value2 = 0.9991044218282156
$ws.Cells.Item(29, 20).Value2 = 0.9999209844048892
$ws.Cells.Item(29, 21).Value2 = 0.9997909019662341
$ws.Cells.Item(29, 22).Value2 = 1.000202790049045
$ws.Cells.Item(29, 23).Value2 = 0.999815817398235

$ws.Cells.Item(30, 3).Value2 = 0.9972430042087267
$ws.Cells.Item(30, 4).Value2 = 1.000724560131472
$ws.Cells.Item(30, 5).Value2 = 1.001775217940601
$ws.Cells.Item(30, 6).Value2 = 0.9994712029435705
$ws.Cells.Item(30, 7).Value2 = 0.9972430042087267
$ws.Cells.Item(30, 8).Value2 = 1.000669636801303
$ws.Cells.Item(30, 9).Value2 = 0.998213508265255
$ws.Cells.Item(30, 10).Value2 = 1.001775217940601
$ws.Cells.Item(30, 11).Value2 = 1.001775217940601
$ws.Cells.Item(30, 12).Value2 = 1.0010270844067
$ws.Cells.Item(30, 13).Value2 = 0.9994097142742324
$ws.Cells.Item(30, 14).Value2 = 1.001775217940601
$ws.Cells.Item(30, 15).Value2 = 1.000724560131472
$ws.Cells.Item(30, 16).Value2 = 0.9989837821700993
$ws.Cells.Item(30, 17).Value2 = 1.000067137202852
$ws.Cells.Item(30, 18).Value2 = 0.9999142607602666
$ws.Cells.Item(30, 19).Value2 = 0.9991257595381438
$ws.Cells.Item(30, 20).Value2 = 0.9999142607602666
$ws.Cells.Item(30, 21).Value2 = 0.9997881241387581
$ws.Cells.Item(30, 22).Value2 = 1.000185542899127
$ws.Cells.Item(30, 23).Value2 = 0.9998167411214824

$ws.Cells.Item(31, 3).Value2 = 1.005525112434104
$ws.Cells.Item(31, 4).Value2 = 0.998597375154899
$ws.Cells.Item(31, 5).Value2 = 0.9963545432700812
$ws.Cells.Item(31, 6).Value2 = 1.001091370237678
$ws.Cells.Item(31, 7).Value2 = 1.005525112434104
$ws.Cells.Item(31, 8).Value2 = 0.9986305896284924
$ws.Cells.Item(31, 9).Value2 = 1.003588777315701
$ws.Cells.Item(31, 10).Value2 = 0.9963545432700812
$ws.Cells.Item(31, 11).Value2 = 0.9963545432700812
$ws.Cells.Item(31, 12).Value2 = 0.9979468989794275
$ws.Cells.Item(31, 13).Value2 = 1.001185932650153
$ws.Cells.Item(31, 14).Value2 = 0.9963545432700812
$ws.Cells.Item(31, 15).Value2 = 0.998597375154899
$ws.Cells.Item(31, 16).Value2 = 1.002061243794501
$ws.Cells.Item(31, 17).Value2 = 0.999891653902526
$ws.Cells.Item(31, 18).Value2 = 1.000159010286361
$ws.Cells.Item(31, 19).Value2 = 1.001769473413052
$ws.Cells.Item(31, 20).Value2 = 1.000159010286361
$ws.Cells.Item(31, 21).Value2 = 1.000415740877309
$ws.Cells.Item(31, 22).Value2 = 0.9996035013558636
$ws.Cells.Item(31, 23).Value2 = 1.000365074958817

